# Append new game-session rows (34-42) to the "data" worksheet, and update
# the sheet view (scroll position / active selection) to match the final
# state of the workbook after the edit.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("data")

$ws.Range("A34").Value = 2
$ws.Range("B34").Value = 'Fall of the Hulks'
$ws.Range("C34").Value = 'Thanos'
$ws.Range("D34").Value = 'Infinity Gems|Emissaries of Evil'
$ws.Range("E34").Value = 'Sapien League'
$ws.Range("F34").Value = 'Nul, Breaker of Worlds (FI)|Agent X-13 (C75)|Totally Awesome Hulk (CH)|Namora (WW)|Iron Fist (DC)'
$ws.Range("G34").Value = 1
$ws.Range("H34").Value = '8|43'
$ws.Range("I34").Value = 'yes'
$ws.Range("K34").Value = 'One card left in villain deck, 1 wound left in stack. Reality gem is weak. One player managed to hit all the henchmen for thinning.'

$ws.Range("A35").Value = 2
$ws.Range("B35").Value = 'Secret Empire of Betrayal'
$ws.Range("C35").Value = 'Fin Fang Foom'
$ws.Range("D35").Value = 'Monsters Unleashed|Shi''ar Imperial Guard'
$ws.Range("E35").Value = 'Shi''ar Patrol Craft'
$ws.Range("F35").Value = 'Sunspot (NM)|Longshot (XM)|Solo (DP)|Amadeus Cho (WW)|Rogue (B)'
$ws.Range("G35").Value = 0
$ws.Range("H35").Value = '34|41'
$ws.Range("I35").Value = 'yes'
$ws.Range("J35").Value = 'Cable as extra hero'
$ws.Range("K35").Value = '10th twist resolved relatively early with still 15+ cards in villain deck.'

$ws.Range("A36").Value = 2
$ws.Range("B36").Value = 'Save Humanity'
$ws.Range("C36").Value = 'The Hood'
$ws.Range("D36").Value = 'Hood''s Gang|S.H.I.E.L.D. Elite'
$ws.Range("E36").Value = 'Mandarin''s Rings'
$ws.Range("F36").Value = 'Cloak & Dagger (CW)|Phoenix Force Cyclops (SW2)|Valkyrie (HOA)|Tigra (CW)|Captain Marvel (SW1)'
$ws.Range("G36").Value = 1
$ws.Range("H36").Value = '38|34'
$ws.Range("I36").Value = 'no'
$ws.Range("K36").Value = 'Few bystanders appeared in the hq and there were quite a number of yellow heroes.'

$ws.Range("A37").Value = 2
$ws.Range("B37").Value = 'Secret Invasion of the Skrull Shapeshifters'
$ws.Range("C37").Value = 'M.O.D.O.K.'
$ws.Range("D37").Value = 'Skrulls|Intelligencia'
$ws.Range("E37").Value = 'Mandroid'
$ws.Range("F37").Value = 'Photon (R)|Black Knight (AM)|Howard the Duck (D)|Maximus (SW1)|Venomized Dr. Strange (VE)|Karma (NM)'
$ws.Range("G37").Value = 1
$ws.Range("H37").Value = '39|31'
$ws.Range("I37").Value = 'yes'
$ws.Range("K37").Value = 'Difficult to consistently get to 8 recruit after all strikes are gone'

$ws.Range("A38").Value = 2
$ws.Range("B38").Value = 'Deadpool Kills the Marvel Universe'
$ws.Range("C38").Value = 'Hela, Goddess of Death'
$ws.Range("D38").Value = 'Omens of Ragnarok|Manhattan (Earth-1610)'
$ws.Range("E38").Value = 'Maggia Goons'
$ws.Range("F38").Value = 'Stingray (DP)|Nova (CH)|Agent Venom (SW2)|Angel Noir (N)'
$ws.Range("G38").Value = 0
$ws.Range("H38").Value = '8|33'
$ws.Range("I38").Value = 'yes'
$ws.Range("K38").Value = 'Final twist one turn too soon. Played a few turns already with an empty hero deck.'

$ws.Range("A39").Value = 2
$ws.Range("B39").Value = 'The Dark World of Svartalfheim'
$ws.Range("C39").Value = 'Immortal Emperor Zheng Zhu'
$ws.Range("D39").Value = 'K''un-Lun|Manhattan (Earth-1610)'
$ws.Range("E39").Value = 'Savage Land Mutates'
$ws.Range("F39").Value = 'Blade (DC)|Dazzler (XM)|Colossus & Wolverine (XM)|Captain Marvel (SW2)|Venom (V)'
$ws.Range("G39").Value = 1
$ws.Range("H39").Value = '21|22'
$ws.Range("I39").Value = 'not really'
$ws.Range("J39").Value = 'First game using Tabletop Simulator'
$ws.Range("K39").Value = 'One twist out only, but it was turns away.'

$ws.Range("A40").Value = 2
$ws.Range("B40").Value = 'Crush Them With My Bare Hands'
$ws.Range("C40").Value = 'Supreme Intelligence of the Kree'
$ws.Range("D40").Value = 'Kree Starforce|Lethal Legion'
$ws.Range("E40").Value = 'Hellfire Cult'
$ws.Range("F40").Value = 'Mr. Fantastic (FF)|Rocket Raccoon (GG)|Beast (XM)|Howard the Duck (D)|Superior Iron Man (SW1)'
$ws.Range("G40").Value = 1
$ws.Range("H40").Value = '39|53'
$ws.Range("I40").Value = 'yes'
$ws.Range("K40").Value = 'The mastermind managed to gather a lot of shards. It took Mr Fantastic''s ultimate and focus accumulation for a few big hits in the late game.'

$ws.Range("A41").Value = 2
$ws.Range("B41").Value = 'Capture Baby Hope'
$ws.Range("C41").Value = 'Morgan Le Fay'
$ws.Range("D41").Value = 'Queen''s Vengeance|Goblin''s Freak Show'
$ws.Range("E41").Value = 'Shi''ar Patrol Craft'
$ws.Range("F41").Value = 'Kitty Pryde (XM)|Corvus Glaive (SW2)|Hulk (B)|Viv Vision (CH)|Wasp (AM)'
$ws.Range("G41").Value = 1
$ws.Range("H41").Value = '46|24'
$ws.Range("I41").Value = 'no'
$ws.Range("K41").Value = 'Lots of ultimate heroes early on. Two twists stacked, but with only two twists out and city mostly empty.'

$ws.Range("A42").Value = 2
$ws.Range("B42").Value = 'The God-Emperor of Battleworld'
$ws.Range("C42").Value = 'Hybrid'
$ws.Range("D42").Value = 'Life Foundation|Hellfire Club'
$ws.Range("E42").Value = 'Doombot Legion'
$ws.Range("F42").Value = 'Gamora (GG)|Karma (NM)|The Warriors Three (HOA)|Ms. Marvel (CH)|Nul, Breaker of Worlds (FI)'
$ws.Range("G42").Value = 1
$ws.Range("H42").Value = '70|36'
$ws.Range("I42").Value = 'no'
$ws.Range("K42").Value = 'Two early master strikes were duds, the emperor was at 11 for a long time and got beaten before any more twists. Nul is really nasty.'


# Update the sheet view to match the saved workbook state: scrolled so
# row 19 is the top-left visible row, with the active selection on H43
# (one cell below/right of the last data row, matching Excel's behavior
# of leaving the selection just past the newly entered data).
$ws.Activate()
$excel.ActiveWindow.SetTopLeftVisibleCell($ws.Range("A19"))
$ws.Range("H43").Select()
